# Update core ABM pharmacokinetics parameters: k_elim (col C) and c_max (col D)
# for osimertinib (row 2) and gefitinib (row 3).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2
$ws.Range("D2").Value = 100
$ws.Range("C3").Value = 0.5
$ws.Range("D3").Value = 100

$ws.Range("D3").Select()
